$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> [open_price(D), close_price(E), high_price(F), low_price(G)]
# shares_outstanding (H) and fixed_ticker (I) are reset to FTNT's own values (766266033 / "FTNT")
# for every data row, undoing cross-contamination from other tickers' rows.
$rowData = @{
    2 = @(6.99399995803833, 7.547999858856201, 7.797999858856201, 6.679999828338623)
    3 = @(8.368000030517578, 9.548000335693359, 9.718000411987305, 7.808000087738037)
    4 = @(8.498000144958496, 6.872000217437744, 8.932000160217285, 6.581999778747559)
    5 = @(6.107999801635742, 5.627999782562256, 6.197999954223633, 4.728000164031982)
    6 = @(6.050000190734863, 6.501999855041504, 6.868000030517578, 5.616000175476074)
    7 = @(6.300000190734863, 6.938000202178955, 7.461999893188477, 6.184000015258789)
    8 = @(7.380000114440918, 6.412000179290772, 7.495999813079834, 5.800000190734863)
    9 = @(6.111999988555908, 6.651999950408936, 6.677999973297119, 5.980000019073486)
    10 = @(7.710000038146973, 7.800000190734863, 8.128000259399414, 7.386000156402588)
    11 = @(7.533999919891357, 7.381999969482422, 8.312000274658203, 7.360000133514404)
    12 = @(7.182000160217285, 7.881999969482422, 8.112000465393066, 7.182000160217285)
    14 = @(10.62199974060059, 11.07199954986572, 11.5959997177124, 10.4060001373291)
    15 = @(12.43400001525879, 12.58199977874756, 13.76599979400635, 12.38000011444092)
    16 = @(18.60000038146973, 16.43600082397461, 18.87400054931641, 15.20400047302246)
    17 = @(13.71599960327148, 15.31400012969971, 15.51399993896484, 13.01200008392334)
    18 = @(16.91799926757812, 18.68400001525879, 19.39200019836425, 16.48800086975098)
    19 = @(15.62199974060059, 16.0620002746582, 17.72999954223633, 15.37600040435791)
    20 = @(15.51599979400635, 16.3120002746582, 16.69400024414062, 14.90999984741211)
    21 = @(21.52000045776367, 23.07200050354004, 24.08799934387207, 21.47599983215332)
    22 = @(19.38800048828125, 21.54800033569336, 22.99600028991699, 19.20199966430664)
    23 = @(27.45400047302246, 27.65999984741211, 30.38999938964844, 24.8479995727539)
    24 = @(23.72800064086914, 22.07399940490723, 26.98800086975098, 21.35000038146973)
    26 = @(37.35400009155273, 40.84600067138672, 41.87799835205078, 37.09400177001953)
    27 = @(47.57199859619141, 54.44800186157227, 55.40000152587891, 47.27199935913086)
    28 = @(58.23400115966797, 67.26799774169922, 69.05000305175781, 56.95399856567383)
    29 = @(72.42400360107422, 59.44800186157227, 72.42400360107422, 52.56800079345703)
    30 = @(68.20200347900391, 57.80199813842773, 71.51599884033203, 57.52799987792969)
    31 = @(56.40000152587891, 59.65000152587891, 63.56000137329102, 55.20999908447266)
    32 = @(49.58000183105469, 57.15999984741211, 57.34000015258789, 45.7400016784668)
    33 = @(49.52999877929688, 52.34000015258789, 54.56999969482422, 45.88000106811523)
    34 = @(65.81999969482422, 63.04999923706055, 69.06999969482422, 61.72999954223633)
    35 = @(75.30000305175781, 77.72000122070312, 81.23999786376953, 72.44999694824219)
    36 = @(58.79000091552734, 57.16999816894531, 60.59999847412109, 55.45999908447266)
    37 = @(58.09999847412109, 64.48999786376953, 67.36000061035156, 56.93999862670898)
    38 = @(68.73000335693359, 63.18000030517578, 73.62999725341797, 62.61000061035156)
    39 = @(60.43000030517578, 58.04000091552734, 61.47999954223633, 56.41999816894531)
    40 = @(78.20999908447266, 78.66000366210938, 83.76999664306641, 76)
    41 = @(95.36000061035156, 100.879997253418, 102.870002746582, 90.8000030517578)
    42 = @(95.22000122070312, 103.7600021362305, 104.5400009155273, 81.69999694824219)
    43 = @(105.3600006103516, 99.90000152587891, 108.7699966430664, 98.6999969482422)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
    $ws.Range("F$row").Value = $vals[2]
    $ws.Range("G$row").Value = $vals[3]
    $ws.Range("H$row").Value = 766266033
    $ws.Range("I$row").Value = "FTNT"
}
